# Updated progress and middle east v0.5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New character rows (15-23), written column-major within each row
#     (A, B, C, E, [G], H, T) so new shared strings are allocated in the
#     same order as the target workbook. ---

# Row 15 - Jahangir bin Ali (Ak Koyunlu)
$ws.Range("A15").Value = "CKU_ak_koyunlu_1"
$ws.Range("B15").Value = "(31, 49, 138)"
$ws.Range("C15").Value = "Jahangir bin Ali"
$ws.Range("E15").Value = "1414.1.1"
$ws.Range("H15").Value = "No"
$ws.Range("T15").Value = "dynn_ak_koyunlu"

# Row 16 - Ioannes (Trebizond)
$ws.Range("A16").Value = "CKU_trebizond_1"
$ws.Range("B16").Value = "(40, 164, 157)"
$ws.Range("C16").Value = "Ioannes"
$ws.Range("E16").Value = "1403.1.1"
$ws.Range("H16").Value = "No"
$ws.Range("T16").Value = "dynn_komnenos"

# Row 17 - Isfahan (Kara Koyunlu)
$ws.Range("A17").Value = "CKU_kara_koyunlu_1"
$ws.Range("B17").Value = "(54, 117, 136)"
$ws.Range("C17").Value = "Isfahan"
$ws.Range("E17").Value = "1410.1.1"
$ws.Range("H17").Value = "No"
$ws.Range("T17").Value = "dynn_kara_koyunlu"

# Row 18 - Bidlo (Bitlis)
$ws.Range("A18").Value = "CKU_bitlis_1"
$ws.Range("B18").Value = "(205, 110, 47)"
$ws.Range("C18").Value = "Bidlo"
$ws.Range("E18").Value = "1421.1.1"
$ws.Range("G18").Value = "(54, 117, 136)"
$ws.Range("H18").Value = "No"
$ws.Range("T18").Value = "dynn_rojaki"

# Row 19 - Albast (Kharabakh)
$ws.Range("A19").Value = "CKU_kharabakh_1"
$ws.Range("B19").Value = "(117, 167, 117)"
$ws.Range("C19").Value = "Albast"
$ws.Range("E19").Value = "1425.1.1 "
$ws.Range("G19").Value = "(54, 117, 136)"
$ws.Range("H19").Value = "No"
$ws.Range("T19").Value = "dynn_hasan_jalalyan"

# Row 20 - Asad al-din Zarin (Hakkari) - name cell carries a custom font
$ws.Range("A20").Value = "CKU_hakkari_1"
$ws.Range("B20").Value = "(107, 159, 136)"
$ws.Range("C20").Value = "Asad al-din Zarin"
$ws.Range("E20").Value = "1412.1.1"
$ws.Range("G20").Value = "(54, 117, 136)"
$ws.Range("H20").Value = "No"
$ws.Range("T20").Value = "dynn_cang"

# Row 21 - Hussain bin Ala-ud-Daulah Ahmed (Jalair)
$ws.Range("A21").Value = "CKU_jalair_1"
$ws.Range("B21").Value = "(162, 165, 193)"
$ws.Range("C21").Value = "Hussain bin Ala-ud-Daulah Ahmed"
$ws.Range("E21").Value = "1407.1.1"
$ws.Range("G21").Value = "(54, 117, 136)"
$ws.Range("H21").Value = "No"
$ws.Range("T21").Value = "dynn_jalayir"

# Row 22 - Muhammed (Mushasha)
$ws.Range("A22").Value = "CKU_mushasha_1"
$ws.Range("B22").Value = "(140, 102, 152)"
$ws.Range("C22").Value = "Muhammed"
$ws.Range("E22").Value = "1422.1.1"
$ws.Range("G22").Value = "(54, 117, 136)"
$ws.Range("H22").Value = "No"
$ws.Range("T22").Value = "dynn_falah"

# Row 23 - Salih Salah ad-Din (Hisn Kayfa)
$ws.Range("A23").Value = "CKU_hisn_kayfa_1"
$ws.Range("B23").Value = "(234, 179, 30)"
$ws.Range("C23").Value = "Salih Salah ad-Din"
$ws.Range("E23").Value = "1411.1.1"
$ws.Range("H23").Value = "No"
$ws.Range("T23").Value = "dynn_ayyubid"

# --- Distinct font colour on the Hakkari ruler's name (new font + cellXfs entry) ---
$ws.Range("C20").Font.Color = 2236704   # BGR packed value of RGB(0x20,0x21,0x22) -> FF202122

# --- Column width tweaks (re-measured / widened columns) ---
# ColumnWidth setter adds a fixed 5/6 offset internally, so back that out of
# the desired stored <col width="..."> value.
$offset = 0.8333333333333334
$ws.Columns.Item(1).ColumnWidth = 18.5546875 - $offset
$ws.Columns.Item(2).ColumnWidth = 13.5546875 - $offset
$ws.Columns.Item(3).ColumnWidth = 28.6640625 - $offset
$ws.Columns.Item(4).ColumnWidth = 4.33203125 - $offset
$ws.Columns.Item(5).ColumnWidth = 9.88671875 - $offset
$ws.Columns.Item(6).ColumnWidth = 11.88671875 - $offset
$ws.Columns.Item(7).ColumnWidth = 11.44140625 - $offset
$ws.Columns.Item(8).ColumnWidth = 7.33203125 - $offset
$ws.Columns.Item(9).ColumnWidth = 7.33203125 - $offset
$ws.Columns.Item(10).ColumnWidth = 8.44140625 - $offset
$ws.Columns.Item(11).ColumnWidth = 10.109375 - $offset
$ws.Columns.Item(12).ColumnWidth = 8 - $offset
$ws.Columns.Item(13).ColumnWidth = 11.88671875 - $offset
$ws.Columns.Item(14).ColumnWidth = 8.33203125 - $offset
$ws.Columns.Item(15).ColumnWidth = 6.44140625 - $offset
$ws.Columns.Item(16).ColumnWidth = 7.5546875 - $offset
$ws.Columns.Item(17).ColumnWidth = 22.33203125 - $offset
$ws.Columns.Item(18).ColumnWidth = 5.109375 - $offset
$ws.Columns.Item(19).ColumnWidth = 7.33203125 - $offset
$ws.Columns.Item(20).ColumnWidth = 7.88671875 - $offset
$ws.Columns.Item(21).ColumnWidth = 14.44140625 - $offset
$ws.Columns.Item(22).ColumnWidth = 14.44140625 - $offset
$ws.Columns.Item(23).ColumnWidth = 9 - $offset
$ws.Columns.Item(24).ColumnWidth = 6.6640625 - $offset
$ws.Columns.Item(25).ColumnWidth = 7.6640625 - $offset
$ws.Columns.Item(26).ColumnWidth = 5.6640625 - $offset
$ws.Columns.Item(27).ColumnWidth = 5.6640625 - $offset

# --- Selection moved to G26 (matches the saved cursor position in the workbook) ---
$ws.Range("G26").Select()
